# Canary tutorial deck edits (quickstart-iter8-process.pptx, slide 1):
#   1. Split the "Query metrics from New Relic and Prometheus" run into
#      three runs: "Query " / "metrics from New Relic and " / "Prometheus".
#   2. Collapse the triple space in
#      "Iter8 experiment   with A/B testing and progressive deployment"
#      down to a single space.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the two target shapes by name so the script isn't dependent on a
# particular z-order / Shapes.Item index.
$queryShape = $null
$iter8Shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 26") { $queryShape = $shp }
    if ($shp.Name -eq "Rectangle 29") { $iter8Shape = $shp }
}

# --- Change 1: Rectangle 26 ("Query metrics from New Relic and Prometheus") ---
$queryRange = $queryShape.TextFrame.TextRange

$part1 = "Query "
$part2 = "metrics from New Relic and "
$part3 = "Prometheus"

$run1 = $queryRange.Characters(1, $part1.Length)
$run1.Text = $part1

$run2 = $queryRange.Characters($part1.Length + 1, $part2.Length)
$run2.Text = $part2

$run3 = $queryRange.Characters($part1.Length + $part2.Length + 1, $part3.Length)
$run3.Text = $part3

# --- Change 2: Rectangle 29 ("Iter8 experiment   with A/B testing and progressive deployment") ---
$iter8Range = $iter8Shape.TextFrame.TextRange

$oldFragment = "Iter8 experiment   with A/B testing and progressive deployment"
$newFragment = "Iter8 experiment with A/B testing and progressive deployment"
$fragIndex = $iter8Range.Text.IndexOf($oldFragment)

$fragRange = $iter8Range.Characters($fragIndex + 1, $oldFragment.Length)
$fragRange.Text = $newFragment
